$wb = $excel.ActiveWorkbook

# Sheet 1
$ws = $wb.Worksheets.Item(1)
$ws.Range("F2").Value = 50
$ws.Range("F7").Value = 5601
$ws.Range("F8").Value = 7568
$ws.Range("F9").Value = 9
$ws.Range("F12").Value = 3823
$ws.Range("F13").Value = 70
$ws.Range("F16").Value = 195
$ws.Range("F21").Value = 592
$ws.Range("F22").Value = 3869
$ws.Range("F23").Value = 129
$ws.Range("F25").Value = 5258
$ws.Range("F27").Value = 2082
$ws.Range("F28").Value = 130
$ws.Range("F30").Value = 7798
$ws.Range("F34").Value = 2170
$ws.Range("F36").Value = 1189
$ws.Range("F44").Value = 29
$ws.Range("F45").Value = 1324
$ws.Range("F46").Value = 2052
$ws.Range("F47").Value = 123
$ws.Range("F48").Value = 217

# Sheet 2
$ws = $wb.Worksheets.Item(2)
$ws.Range("F11").Value = 121

# Sheet 3
$ws = $wb.Worksheets.Item(3)
$ws.Range("F2").Value = 563
$ws.Range("F3").Value = 736

# Sheet 4
$ws = $wb.Worksheets.Item(4)
$ws.Range("F2").Value = 50
$ws.Range("F5").Value = 563
$ws.Range("F6").Value = 736
$ws.Range("F8").Value = 5601
$ws.Range("F9").Value = 7568
$ws.Range("F10").Value = 9
$ws.Range("F11").Value = 3823
$ws.Range("F14").Value = 195
$ws.Range("F20").Value = 592
$ws.Range("F21").Value = 3869
$ws.Range("F23").Value = 129
$ws.Range("F25").Value = 5258
$ws.Range("F27").Value = 2082
$ws.Range("F28").Value = 130
$ws.Range("F30").Value = 7798
$ws.Range("F34").Value = 2170
$ws.Range("F36").Value = 1189
$ws.Range("F42").Value = 29
$ws.Range("F43").Value = 1324
$ws.Range("F44").Value = 2052
$ws.Range("F45").Value = 123
$ws.Range("F47").Value = 217
